# Update odds values on row 2 (Brisbane Roar vs Melbourne Victory)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.62
$ws.Range("G2").Value = 2.64
$ws.Range("I2").Value = 2.96
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 3.65
$ws.Range("V2").Value = 1.51
$ws.Range("W2").Value = 1.6
$ws.Range("Y2").Value = 14
$ws.Range("AF2").Value = 22
$ws.Range("AG2").Value = 13.5
$ws.Range("AI2").Value = 55
$ws.Range("AJ2").Value = 48
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 48
$ws.Range("AN2").Value = 28
$ws.Range("AO2").Value = 34

# Update odds values on row 3 (Mgladbach vs RB Leipzig)
$ws.Range("H3").Value = 2.22
$ws.Range("AD3").Value = 13
$ws.Range("AG3").Value = 17
$ws.Range("AN3").Value = 18

# Update odds values on row 7 (Como vs Sassuolo)
$ws.Range("K7").Value = 4.2
$ws.Range("P7").Value = 1.94
$ws.Range("Q7").Value = 2.02
$ws.Range("AF7").Value = 9.6
$ws.Range("AJ7").Value = 19
$ws.Range("AK7").Value = 23

# Update odds values on row 8 (Getafe vs Elche)
$ws.Range("G8").Value = 2.36
$ws.Range("O8").Value = 1.65
